$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 148-149),
# shifting the existing rows 148:226 down to 150:228.
$ws.Rows("148:149").Insert()

# Row 148 - new weekly record (Primera)
$ws.Cells.Item(148, 1).Value  = 9
$ws.Cells.Item(148, 2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(148, 3).Value  = 'Metropolitana'
$ws.Cells.Item(148, 4).Value  = 44455
$ws.Cells.Item(148, 5).Value  = 13
$ws.Cells.Item(148, 6).Value  = 100112012
$ws.Cells.Item(148, 7).Value  = 'Espinaca'
$ws.Cells.Item(148, 8).Value  = 'Sin especificar'
$ws.Cells.Item(148, 9).Value  = 'Primera'
$ws.Cells.Item(148, 10).Value = 250
$ws.Cells.Item(148, 11).Value = 6000
$ws.Cells.Item(148, 12).Value = 7000
$ws.Cells.Item(148, 13).Value = 6500
$ws.Cells.Item(148, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(148, 15).Value = 'Provincia de Chacabuco'
$ws.Cells.Item(148, 16).Value = 650
$ws.Cells.Item(148, 17).Value = 10
$ws.Cells.Item(148, 18).Value = 'Hortaliza'

# Row 149 - new weekly record (Segunda)
$ws.Cells.Item(149, 1).Value  = 9
$ws.Cells.Item(149, 2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(149, 3).Value  = 'Metropolitana'
$ws.Cells.Item(149, 4).Value  = 44455
$ws.Cells.Item(149, 5).Value  = 13
$ws.Cells.Item(149, 6).Value  = 100112012
$ws.Cells.Item(149, 7).Value  = 'Espinaca'
$ws.Cells.Item(149, 8).Value  = 'Sin especificar'
$ws.Cells.Item(149, 9).Value  = 'Segunda'
$ws.Cells.Item(149, 10).Value = 106
$ws.Cells.Item(149, 11).Value = 4000
$ws.Cells.Item(149, 12).Value = 5000
$ws.Cells.Item(149, 13).Value = 4500
$ws.Cells.Item(149, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(149, 15).Value = 'Provincia de Chacabuco'
$ws.Cells.Item(149, 16).Value = 450
$ws.Cells.Item(149, 17).Value = 10
$ws.Cells.Item(149, 18).Value = 'Hortaliza'
